# Update "Clan Games" data - 2025-12-22
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("clan games")

# Wrap every average formula in column C (rows 2-49) with ROUND(..., 0)
for ($r = 2; $r -le 49; $r++) {
    $ws.Range("C$r").Formula = "=ROUND(AVERAGE(F$r`:AB$r), 0)"
}

# Update the 22/12/2025 (column H) contributions with the new values
$hUpdates = @{
    6  = 900
    11 = 3000
    12 = 6650
    15 = 10100
    19 = 4450
    24 = 900
    26 = 4200
    28 = 1900
    31 = 20000
    34 = 10000
    38 = 1200
    39 = 800
}

foreach ($row in $hUpdates.Keys) {
    $ws.Range("H$row").Value = $hUpdates[$row]
}
